$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 3335  # was 3318
$ws.Range("F5").Value = 212  # was 211
$ws.Range("F6").Value = 4824  # was 4812
$ws.Range("F7").Value = 466  # was 465
$ws.Range("F8").Value = 297  # was 293
$ws.Range("F10").Value = 627  # was 625
$ws.Range("F12").Value = 29  # was 26
$ws.Range("F13").Value = 12  # was 11
$ws.Range("F14").Value = 657  # was 656
$ws.Range("F15").Value = 287  # was 286
$ws.Range("F16").Value = 25  # was 23
$ws.Range("F18").Value = 143  # was 142
$ws.Range("F20").Value = 4753  # was 4747
$ws.Range("F21").Value = 20  # was 15
$ws.Range("F22").Value = 34  # was 33
$ws.Range("F24").Value = 5893  # was 5885
$ws.Range("F27").Value = 239  # was 237
$ws.Range("F28").Value = 668  # was 666
$ws.Range("F29").Value = 4419  # was 4418
$ws.Range("F31").Value = 94  # was 93
$ws.Range("F33").Value = 857  # was 855
$ws.Range("F35").Value = 6  # was 3
$ws.Range("F36").Value = 784  # was 783
$ws.Range("F37").Value = 833  # was 826

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 13  # was 12

# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 38  # was 37

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 38  # was 37
$ws.Range("F8").Value = 3335  # was 3318
$ws.Range("F9").Value = 212  # was 211
$ws.Range("F10").Value = 4824  # was 4812
$ws.Range("F11").Value = 466  # was 465
$ws.Range("F12").Value = 297  # was 293
$ws.Range("F14").Value = 627  # was 625
$ws.Range("F16").Value = 29  # was 26
$ws.Range("F17").Value = 12  # was 11
$ws.Range("F18").Value = 657  # was 656
$ws.Range("F19").Value = 287  # was 286
$ws.Range("F20").Value = 25  # was 23
$ws.Range("F23").Value = 143  # was 142
$ws.Range("F25").Value = 4753  # was 4747
$ws.Range("F26").Value = 20  # was 15
$ws.Range("F27").Value = 34  # was 33
$ws.Range("F29").Value = 5893  # was 5885
$ws.Range("F32").Value = 239  # was 237
$ws.Range("F33").Value = 668  # was 666
$ws.Range("F34").Value = 4419  # was 4418
$ws.Range("F36").Value = 13  # was 12
$ws.Range("F37").Value = 94  # was 93
$ws.Range("F39").Value = 857  # was 855
$ws.Range("F41").Value = 6  # was 3
$ws.Range("F42").Value = 784  # was 783
$ws.Range("F43").Value = 833  # was 826
